$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("C3").Value = "18.01: Kontakttime, kursansvarlig tilgjengelig i Aud Max og på [Zoom](https://nhh.zoom.us/j/66065667678?pwd=ME1LK294VUw4SEt3eHI2V1ZuZm5MZz09)."
$ws.Range("D3").Value = "20.01: Oppgaveseminar i Aud Max og på [Zoom](https://nhh.zoom.us/j/66065667678?pwd=ME1LK294VUw4SEt3eHI2V1ZuZm5MZz09). Se \@ref(seminar) for oppgaver."
$ws.Range("C4").Value = "25.01: Kontakttime, kursansvarlig tilgjengelig i Aud Max og på [Zoom](https://nhh.zoom.us/j/66065667678?pwd=ME1LK294VUw4SEt3eHI2V1ZuZm5MZz09)."
$ws.Range("D4").Value = "27.01: Kontakttime, kursansvarlig tilgjengelig i Aud Max og på [Zoom](https://nhh.zoom.us/j/66065667678?pwd=ME1LK294VUw4SEt3eHI2V1ZuZm5MZz09)."
$ws.Range("C5").Value = "01.02: Kontakttime, kursansvarlig tilgjengelig i Aud Max og på [Zoom](https://nhh.zoom.us/j/66065667678?pwd=ME1LK294VUw4SEt3eHI2V1ZuZm5MZz09)."
$ws.Range("D5").Value = "03.02: Oppgaveseminar i Aud Max og på [Zoom](https://nhh.zoom.us/j/66065667678?pwd=ME1LK294VUw4SEt3eHI2V1ZuZm5MZz09). Se \@ref(seminar) for oppgaver."
$ws.Range("C6").Value = "08.02: Kontakttime, kursansvarlig tilgjengelig i Aud Max og på [Zoom](https://nhh.zoom.us/j/66065667678?pwd=ME1LK294VUw4SEt3eHI2V1ZuZm5MZz09)."
$ws.Range("D6").Value = "10.02: **Oversiktsforelesning: Hypotesetesting** i Aud Max og på [Zoom](https://nhh.zoom.us/j/66065667678?pwd=ME1LK294VUw4SEt3eHI2V1ZuZm5MZz09)."
$ws.Range("C7").Value = "15.02: Kontakttime, kursansvarlig tilgjengelig i Aud Max og på [Zoom](https://nhh.zoom.us/j/66065667678?pwd=ME1LK294VUw4SEt3eHI2V1ZuZm5MZz09)."
$ws.Range("D7").Value = "17.02: Oppgaveseminar i Aud Max og på [Zoom](https://nhh.zoom.us/j/66065667678?pwd=ME1LK294VUw4SEt3eHI2V1ZuZm5MZz09). Se \@ref(seminar) for oppgaver."
$ws.Range("C8").Value = "22.02: Kontakttime, kursansvarlig tilgjengelig i Aud Max og på [Zoom](https://nhh.zoom.us/j/66065667678?pwd=ME1LK294VUw4SEt3eHI2V1ZuZm5MZz09)."
$ws.Range("D8").Value = "24.02: Oppgaveseminar i Aud Max og på [Zoom](https://nhh.zoom.us/j/66065667678?pwd=ME1LK294VUw4SEt3eHI2V1ZuZm5MZz09). Se \@ref(seminar) for oppgaver."
$ws.Range("C9").Value = "01.03: Kontakttime, kursansvarlig tilgjengelig i Aud Max og på [Zoom](https://nhh.zoom.us/j/66065667678?pwd=ME1LK294VUw4SEt3eHI2V1ZuZm5MZz09)."
$ws.Range("D9").Value = "03.03: **Oversiktsforelesning: Regresjon** i Aud Max og på  [Zoom](https://nhh.zoom.us/j/66065667678?pwd=ME1LK294VUw4SEt3eHI2V1ZuZm5MZz09)."
$ws.Range("C10").Value = "08.03: Kontakttime, kursansvarlig tilgjengelig i Aud Max og på [Zoom](https://nhh.zoom.us/j/66065667678?pwd=ME1LK294VUw4SEt3eHI2V1ZuZm5MZz09)."
$ws.Range("D10").Value = "08.03: Oppgaveseminar i Aud Max og på [Zoom](https://nhh.zoom.us/j/66065667678?pwd=ME1LK294VUw4SEt3eHI2V1ZuZm5MZz09). Se \@ref(seminar) for oppgaver."
$ws.Range("C11").Value = "15.03: Kontakttime, kursansvarlig tilgjengelig i Aud Max og på [Zoom](https://nhh.zoom.us/j/66065667678?pwd=ME1LK294VUw4SEt3eHI2V1ZuZm5MZz09)."
$ws.Range("D11").Value = "17.03:  **Oversiktsforelesning: Logistisk, panel, kNN** i Aud Max og på  [Zoom](https://nhh.zoom.us/j/66065667678?pwd=ME1LK294VUw4SEt3eHI2V1ZuZm5MZz09)."
$ws.Range("C12").Value = "22.03: Kontakttime, kursansvarlig tilgjengelig i Aud Max og på [Zoom](https://nhh.zoom.us/j/66065667678?pwd=ME1LK294VUw4SEt3eHI2V1ZuZm5MZz09)."
$ws.Range("D12").Value = "24.03: Oppgaveseminar i Aud Max og på [Zoom](https://nhh.zoom.us/j/66065667678?pwd=ME1LK294VUw4SEt3eHI2V1ZuZm5MZz09). Se \@ref(seminar) for oppgaver."
$ws.Range("C13").Value = "29.03: Kontakttime, kursansvarlig tilgjengelig i Aud Max og på [Zoom](https://nhh.zoom.us/j/66065667678?pwd=ME1LK294VUw4SEt3eHI2V1ZuZm5MZz09)."
$ws.Range("D13").Value = "31.03:  **Oversiktsforelesning: Tidsrekker** i Aud Max og på  [Zoom](https://nhh.zoom.us/j/66065667678?pwd=ME1LK294VUw4SEt3eHI2V1ZuZm5MZz09)."
$ws.Range("C14").Value = "05.04: Kontakttime, kursansvarlig tilgjengelig i Aud Max og på [Zoom](https://nhh.zoom.us/j/66065667678?pwd=ME1LK294VUw4SEt3eHI2V1ZuZm5MZz09)."
$ws.Range("C16").Value = "19.04:  **Oversiktsforelesning: Grunnleggende statistikk/Tips til hj.eksamen** i Aud Max og på  [Zoom](https://nhh.zoom.us/j/66065667678?pwd=ME1LK294VUw4SEt3eHI2V1ZuZm5MZz09)."
$ws.Range("D16").Value = "21.04: **Eksamensoppgaver: Skoleeksamen** i Aud Max og på  [Zoom](https://nhh.zoom.us/j/66065667678?pwd=ME1LK294VUw4SEt3eHI2V1ZuZm5MZz09)."
